# Update the "Förändrad" (Changed) date column C for rows 2 through 28
# from serial date 45424 (2024-05-12) to serial date 45425 (2024-05-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $val = $cell.Value()
    if ($val.ToOADate() -eq 45424) {
        $cell.Value = 45425
    }
}
